$d = $word.ActiveDocument

# Locate the paragraph that ends the "File Handling" section:
# "Reading and writing to text and CSV files"
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd("`r") -eq "Reading and writing to text and CSV files") {
        $target = $p
    }
}

if ($target -eq $null) {
    throw "Could not locate anchor paragraph 'Reading and writing to text and CSV files'"
}

# Insert a fresh empty paragraph right after it, then replace that
# paragraph's content with the full Git/GitHub block via raw OOXML,
# so every run/paragraph gets the exact formatting from the source.
$target.Range.InsertParagraphAfter() | Out-Null
$newPara = $target.Next()
$xmlPayload = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>
  <w:pPr>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1"/>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t>Git</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/>
      <w:b/>
      <w:bCs/>
      <w:sz w:val="24"/>
      <w:lang w:val="es-ES"/>
    </w:rPr>
    <w:t xml:space="preserve"> and GitHub</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>How Git and GitHub work</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
    <w:t>Explanation of the process for including changes in a local repository</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t xml:space="preserve"> (working directory, staging area, and local repository)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>Commands to include changes in a local repository (</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:sz w:val="20"/>
      <w:szCs w:val="20"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>git status, git add, and git commit -m &#8216;comment&#8217;</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>How to push a local repository to a remote repository</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="1"/>
      <w:numId w:val="4"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
      <w:lang w:val="en-US"/>
    </w:rPr>
    <w:t>&#191;</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:cs="Times New Roman"/>
      <w:color w:val="000000" w:themeColor="text1"/>
    </w:rPr>
    <w:t>What is the library called to integrate Git and Python?</w:t>
  </w:r>
</w:p>
<w:p>
  <w:pPr>
    <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1"/>
    <w:ind w:left="1440"/>
  </w:pPr>
</w:p>
</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$newPara.Range.InsertXML($xmlPayload) | Out-Null

"Git and GitHub topics inserted"
